$d = $word.ActiveDocument

# Phase 1: replace each unique original placeholder with a collision-proof temp token.
$targets = @(
  "{0}",
  "{1} {2}",
  "{3}",
  "{4}",
  "{5}",
  "{6} ",
  "{7}",
  "{8}",
  "{9}",
  "{10}",
  "{11}",
  "{12}",
  "{13}",
  "{14}",
  "{15}",
  "{16} ",
  "{17}",
  "{18}",
  "{19}",
  "{20}"
)

for ($i = 0; $i -lt $targets.Count; $i++) {
    $token = "@@TOK$i@@"
    $null = $d.Content.Find.Execute($targets[$i], $true, $false, $false, $false, $false, $true, 1, $false, $token, 2)
}

# Phase 2: replace each temp token with its final content (splitting into multiple runs where needed).
$null = $d.Content.Find.Execute("@@TOK0@@", $true, $false, $false, $false, $false, $true, 1, $false, "{6}", 2)

$rng = $d.Content
$null = $rng.Find.Execute("@@TOK1@@")
$rng.Text = "{7}"
$rng.Collapse(0)
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter(" ")
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter("{8}")

$rng = $d.Content
$null = $rng.Find.Execute("@@TOK2@@")
$rng.Text = "{10}"
$rng.Collapse(0)
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter(" ")
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter("{10'}")

$null = $d.Content.Find.Execute("@@TOK3@@", $true, $false, $false, $false, $false, $true, 1, $false, "{11}", 2)

$null = $d.Content.Find.Execute("@@TOK4@@", $true, $false, $false, $false, $false, $true, 1, $false, "{15}", 2)

$rng = $d.Content
$null = $rng.Find.Execute("@@TOK5@@")
$rng.Text = "{20}"
$rng.Collapse(0)
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter(" ")

$null = $d.Content.Find.Execute("@@TOK6@@", $true, $false, $false, $false, $false, $true, 1, $false, "{21}", 2)

$null = $d.Content.Find.Execute("@@TOK7@@", $true, $false, $false, $false, $false, $true, 1, $false, "{22}", 2)

$null = $d.Content.Find.Execute("@@TOK8@@", $true, $false, $false, $false, $false, $true, 1, $false, "{21}", 2)

$null = $d.Content.Find.Execute("@@TOK9@@", $true, $false, $false, $false, $false, $true, 1, $false, "{10'''}", 2)

$rng = $d.Content
$null = $rng.Find.Execute("@@TOK10@@")
$rng.Text = "{1}"
$rng.Collapse(0)
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter(" ")
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter("{2}")

$rng = $d.Content
$null = $rng.Find.Execute("@@TOK11@@")
$rng.Text = "{"
$rng.Collapse(0)
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter("5")
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter("}")

$null = $d.Content.Find.Execute("@@TOK12@@", $true, $false, $false, $false, $false, $true, 1, $false, "{4}", 2)

$null = $d.Content.Find.Execute("@@TOK13@@", $true, $false, $false, $false, $false, $true, 1, $false, "{5}", 2)

$null = $d.Content.Find.Execute("@@TOK14@@", $true, $false, $false, $false, $false, $true, 1, $false, "{3}", 2)

$rng = $d.Content
$null = $rng.Find.Execute("@@TOK15@@")
$rng.Text = "{3'}"
$rng.Collapse(0)
$rng = $d.Range($rng.End, $rng.End)
$rng.InsertAfter(" ")

$null = $d.Content.Find.Execute("@@TOK16@@", $true, $false, $false, $false, $false, $true, 1, $false, "{5'}", 2)

$null = $d.Content.Find.Execute("@@TOK17@@", $true, $false, $false, $false, $false, $true, 1, $false, "{5''}", 2)

$null = $d.Content.Find.Execute("@@TOK18@@", $true, $false, $false, $false, $false, $true, 1, $false, "{11}", 2)

$null = $d.Content.Find.Execute("@@TOK19@@", $true, $false, $false, $false, $false, $true, 1, $false, "{15}", 2)
